$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2487
$ws1.Range("F7").Value = 1307
$ws1.Range("F8").Value = 1663
$ws1.Range("F13").Value = 149
$ws1.Range("F18").Value = 8458
$ws1.Range("F21").Value = 10529
$ws1.Range("F26").Value = 519
$ws1.Range("F27").Value = 188
$ws1.Range("F30").Value = 18
$ws1.Range("F31").Value = 11
$ws1.Range("F33").Value = 324
$ws1.Range("F34").Value = 421

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2487
$ws4.Range("F10").Value = 1307
$ws4.Range("F12").Value = 1663
$ws4.Range("F18").Value = 149
$ws4.Range("F24").Value = 8458
$ws4.Range("F27").Value = 10529
$ws4.Range("F34").Value = 519
$ws4.Range("F38").Value = 188
$ws4.Range("F47").Value = 421
